$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Rename first sheet, add the two new sheets in order ---
$ws1.Name = "NoState"

$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "State"

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Dependents"

# --- 2. Build the "State" sheet ---
# Header row: reuse columns A-F from NoState (same text + bold/fill style),
# then add the two new "*_State" headers and reuse the "Policy" header.
$ws1.Range("A1:F1").Copy($ws2.Range("A1:F1"))
$ws1.Range("A1").Copy($ws2.Range("G1"))
$ws2.Range("G1").Value = "Home_State"
$ws1.Range("A1").Copy($ws2.Range("H1"))
$ws2.Range("H1").Value = "Host_State"
$ws1.Range("G1").Copy($ws2.Range("I1"))

# Row 2
$ws2.Range("A2").Value = "Test Excel"
$ws2.Range("B2").Value = "Canada"
$ws2.Range("C2").Value = "United States"
$ws2.Range("D2").Value = 87000
$ws2.Range("E2").NumberFormat = "mm-dd-yy"
$ws2.Range("E2").Value = 43101
$ws1.Range("F2").Copy()
$ws2.Range("F2").PasteSpecial(-4163)
$ws2.Range("G2").Value = "Alberta"
$ws2.Range("H2").Value = "Alabama"
$ws2.Range("I2").Value = "CP Plus PPT"

# Row 3
$ws2.Range("A3").Value = "Test Excel 2"
$ws2.Range("B3").Value = "United States"
$ws2.Range("C3").Value = "Canada"
$ws2.Range("D3").Value = 58900
$ws2.Range("E2").Copy($ws2.Range("E3"))
$ws2.Range("E3").Value = 43101
$ws1.Range("F2").Copy()
$ws2.Range("F3").PasteSpecial(-4163)
$ws2.Range("G3").Value = "Alabama"
$ws2.Range("H3").Value = "Alberta"
$ws2.Range("I3").Value = "CP Plus PPT"

$excel.CutCopyMode = $false

# --- 3. Build the "Dependents" sheet ---
# Header row: reuse columns A-F from NoState, reuse "Policy" header, add "Dependents" header.
$ws1.Range("A1:F1").Copy($ws3.Range("A1:F1"))
$ws1.Range("G1").Copy($ws3.Range("G1"))
$ws1.Range("A1").Copy($ws3.Range("H1"))
$ws3.Range("H1").Value = "Dependents"

# Row 2 - identical to NoState row 2 (Scenario..Policy), plus a new Dependents count
$ws1.Range("A2:G2").Copy($ws3.Range("A2:G2"))
$ws3.Range("H2").NumberFormat = "0"
$ws3.Range("H2").Value = 5

# Row 3 - identical to NoState row 3 (Scenario..Policy), plus a new Dependents count
$ws1.Range("A3:G3").Copy($ws3.Range("A3:G3"))
$ws3.Range("H2").Copy($ws3.Range("H3"))
$ws3.Range("H3").Value = 3

$excel.CutCopyMode = $false

# --- 4. Fix up view/selection state ---
$ws1.Range("A1:G6").Select()
$ws2.Range("B2").Select()
$ws3.Range("H2").Select()
